$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 8189.727
$ws.Range("I34").Value = 561
$ws.Range("J34").Value = 28533
$ws.Range("K34").Value = 561
$ws.Range("L34").Value = 28533
$ws.Range("M34").Value = -358
$ws.Range("N34").Value = -28939
$ws.Range("H36").Value = 8189.727
$ws.Range("I36").Value = 561
$ws.Range("J36").Value = 28533
$ws.Range("K36").Value = 561
$ws.Range("L36").Value = 28533
$ws.Range("M36").Value = 154
$ws.Range("N36").Value = -29963
$ws.Range("H40").Value = 2565.8147
$ws.Range("I40").Value = 1962
$ws.Range("J40").Value = 3999.875
$ws.Range("K40").Value = 1962
$ws.Range("L40").Value = 3999.875
$ws.Range("M40").Value = -1787
$ws.Range("N40").Value = -4349.875
$ws.Range("H64").Value = 69172.8
$ws.Range("I64").Value = 102359.2
$ws.Range("J64").Value = 2800
$ws.Range("K64").Value = 102359.2
$ws.Range("L64").Value = 2800
$ws.Range("M64").Value = -102111.2
$ws.Range("N64").Value = -3296
$ws.Range("H67").Value = 69172.8
$ws.Range("I67").Value = 102359.2
$ws.Range("J67").Value = 2800
$ws.Range("K67").Value = 102359.2
$ws.Range("L67").Value = 2800
$ws.Range("M67").Value = -101501.2
$ws.Range("N67").Value = -4516
$ws.Range("H137").Value = 3852.1538
$ws.Range("I137").Value = 1007
$ws.Range("J137").Value = 4294.7334
$ws.Range("K137").Value = 3021
$ws.Range("L137").Value = 12884.2002
$ws.Range("M137").Value = -471
$ws.Range("N137").Value = -17984.2002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29544.373
$ws.Range("I32").Value = 29409.352
$ws.Range("K32").Value = 29409.352
$ws.Range("M32").Value = -29122.352

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 66399.336
$ws.Range("J57").Value = 66399.336
$ws.Range("L57").Value = 66399.336
$ws.Range("N57").Value = -67839.336
$ws.Range("H136").Value = 66399.336
$ws.Range("J136").Value = 66399.336
$ws.Range("L136").Value = 66399.336
$ws.Range("N136").Value = -76599.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 163222.8
$ws.Range("I31").Value = 1864.5
$ws.Range("J31").Value = 217008.89
$ws.Range("K31").Value = 1864.5
$ws.Range("L31").Value = 217008.89
$ws.Range("M31").Value = -1569.5
$ws.Range("N31").Value = -217598.89
$ws.Range("H34").Value = 163222.8
$ws.Range("I34").Value = 1864.5
$ws.Range("J34").Value = 217008.89
$ws.Range("K34").Value = 1864.5
$ws.Range("L34").Value = 217008.89
$ws.Range("M34").Value = -1662.5
$ws.Range("N34").Value = -217412.89
$ws.Range("H58").Value = 1510.091
$ws.Range("I58").Value = 1433.4828
$ws.Range("J58").Value = 1658.2
$ws.Range("K58").Value = 1433.4828
$ws.Range("L58").Value = 1658.2
$ws.Range("M58").Value = -1230.4828
$ws.Range("N58").Value = -2064.2
$ws.Range("H136").Value = 1510.091
$ws.Range("I136").Value = 1433.4828
$ws.Range("J136").Value = 1658.2
$ws.Range("K136").Value = 4300.4484
$ws.Range("L136").Value = 4974.6
$ws.Range("M136").Value = -1750.4484
$ws.Range("N136").Value = -10074.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 600
$ws.Range("J20").Value = 1500
$ws.Range("L20").Value = 4500
$ws.Range("N20").Value = -4954
$ws.Range("H22").Value = 20940
$ws.Range("I22").Value = 900
$ws.Range("J22").Value = 25950
$ws.Range("K22").Value = 2700
$ws.Range("L22").Value = 77850
$ws.Range("M22").Value = -2531
$ws.Range("N22").Value = -78188
$ws.Range("H27").Value = 20940
$ws.Range("I27").Value = 900
$ws.Range("J27").Value = 25950
$ws.Range("K27").Value = 2700
$ws.Range("L27").Value = 77850
$ws.Range("M27").Value = -2598
$ws.Range("N27").Value = -78054
$ws.Range("H39").Value = 1758.5
$ws.Range("J39").Value = 1791.4783
$ws.Range("L39").Value = 5374.4349
$ws.Range("N39").Value = -5962.4349
$ws.Range("H46").Value = 2200
$ws.Range("J46").Value = 2200
$ws.Range("L46").Value = 6600
$ws.Range("N46").Value = -6782
$ws.Range("H58").Value = 1702146.6
$ws.Range("J58").Value = 2042376.2
$ws.Range("L58").Value = 6127128.6
$ws.Range("N58").Value = -6127384.6
$ws.Range("H74").Value = 3439.8
$ws.Range("I74").Value = 1200
$ws.Range("J74").Value = 3999.75
$ws.Range("K74").Value = 3600
$ws.Range("L74").Value = 11999.25
$ws.Range("M74").Value = -2539
$ws.Range("N74").Value = -14121.25
$ws.Range("H77").Value = 3439.8
$ws.Range("I77").Value = 1200
$ws.Range("J77").Value = 3999.75
$ws.Range("K77").Value = 10800
$ws.Range("L77").Value = 35997.75
$ws.Range("M77").Value = -5496
$ws.Range("N77").Value = -46605.75
$ws.Range("H113").Value = 4623.48
$ws.Range("I113").Value = 4942.913
$ws.Range("J113").Value = 950
$ws.Range("K113").Value = 14828.739
$ws.Range("L113").Value = 2850
$ws.Range("M113").Value = -12658.739
$ws.Range("N113").Value = -7190
$ws.Range("H122").Value = 4288.8887
$ws.Range("I122").Value = 381.63635
$ws.Range("J122").Value = 21480.8
$ws.Range("K122").Value = 3434.72715
$ws.Range("L122").Value = 193327.2
$ws.Range("M122").Value = -984.7271499999997
$ws.Range("N122").Value = -198227.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2383.389
$ws.Range("I136").Value = 1764.5714
$ws.Range("J136").Value = 4549.25
$ws.Range("K136").Value = 5293.7142
$ws.Range("L136").Value = 13647.75
$ws.Range("M136").Value = -2743.7142
$ws.Range("N136").Value = -18747.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 80002.39999999999
$ws.Range("J8").Value = 80002.39999999999
$ws.Range("L8").Value = 80002.39999999999
$ws.Range("N8").Value = -80282.39999999999
$ws.Range("H11").Value = 3900
$ws.Range("J11").Value = 3900
$ws.Range("L11").Value = 3900
$ws.Range("N11").Value = -4184
$ws.Range("H108").Value = 26244.666
$ws.Range("J108").Value = 26244.666
$ws.Range("L108").Value = 26244.666
$ws.Range("N108").Value = -33924.666
